# Fruta / hortaliza, semanal
# Inserts two new weekly price-report rows (1087, 1088) into the Lechuga
# (lettuce) price table, pushing the existing rows 1087..1175 down to
# 1089..1177.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 1087, shifting everything
# below (old rows 1087..1175) down to 1089..1177.
$ws.Range("A1087:A1088").EntireRow.Insert()

# --- New row 1087: Lechuga Conconina(o), Primera ---
$ws.Range("A1087").Value2 = 7
$ws.Range("B1087").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C1087").Value2 = "Ñuble"
$ws.Range("D1087").Value2 = 45106
$ws.Range("E1087").Value2 = 16
$ws.Range("F1087").Value2 = 100112033
$ws.Range("G1087").Value2 = "Lechuga"
$ws.Range("H1087").Value2 = "Conconina(o)"
$ws.Range("I1087").Value2 = "Primera"
$ws.Range("J1087").Value2 = 150
$ws.Range("K1087").Value2 = 6000
$ws.Range("L1087").Value2 = 6000
$ws.Range("M1087").Value2 = 6000
$ws.Range("N1087").Value2 = "$/caja 10 unidades"
$ws.Range("O1087").Value2 = "Región del Maule"
$ws.Range("P1087").Value2 = 600
$ws.Range("Q1087").Value2 = 10
$ws.Range("R1087").Value2 = "Hortaliza"

# --- New row 1088: Lechuga Escarola, Primera ---
$ws.Range("A1088").Value2 = 7
$ws.Range("B1088").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C1088").Value2 = "Ñuble"
$ws.Range("D1088").Value2 = 45106
$ws.Range("E1088").Value2 = 16
$ws.Range("F1088").Value2 = 100112033
$ws.Range("G1088").Value2 = "Lechuga"
$ws.Range("H1088").Value2 = "Escarola"
$ws.Range("I1088").Value2 = "Primera"
$ws.Range("J1088").Value2 = 220
$ws.Range("K1088").Value2 = 6000
$ws.Range("L1088").Value2 = 7000
$ws.Range("M1088").Value2 = 6455
$ws.Range("N1088").Value2 = "$/caja 15 unidades"
$ws.Range("O1088").Value2 = "Región del Maule"
$ws.Range("P1088").Value2 = 430
$ws.Range("Q1088").Value2 = 15
$ws.Range("R1088").Value2 = "Hortaliza"
